# Auto-generated edit script applying numeric updates to the Marilith Profits workbook
# Each sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) gets targeted cell updates,
# including a few new cells being populated and a few cells being cleared entirely.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 4312.625
$ws.Range("I2").Value = 1999
$ws.Range("J2").Value = 5083.8335
$ws.Range("K2").Value = 1999
$ws.Range("L2").Value = 5083.8335
$ws.Range("M2").Value = -1886
$ws.Range("N2").Value = -5309.8335
$ws.Range("H4").Value = 133.57143
$ws.Range("I4").Value = 142.5
$ws.Range("J4").Value = 80
$ws.Range("K4").Value = 142.5
$ws.Range("L4").Value = 80
$ws.Range("M4").Value = -28.5
$ws.Range("N4").Value = -308
$ws.Range("H11").Value = 185.625
$ws.Range("I11").Value = 185.625
$ws.Range("K11").Value = 185.625
$ws.Range("M11").Value = -45.625
$ws.Range("H16").Value = 13000
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H80").Value = 1675
$ws.Range("J80").Value = 3000
$ws.Range("L80").Value = 9000
$ws.Range("N80").Value = -10996
$ws.Range("H83").Value = 1675
$ws.Range("J83").Value = 3000
$ws.Range("L83").Value = 27000
$ws.Range("N83").Value = -36984
$ws.Range("H98").Value = 1392.375
$ws.Range("I98").Value = 1161.8572
$ws.Range("J98").Value = 3006
$ws.Range("K98").Value = 1161.8572
$ws.Range("L98").Value = 3006
$ws.Range("M98").Value = 336.1428000000001
$ws.Range("N98").Value = -6002
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H122").Value = 1392.375
$ws.Range("I122").Value = 1161.8572
$ws.Range("J122").Value = 3006
$ws.Range("K122").Value = 3485.5716
$ws.Range("L122").Value = 9018
$ws.Range("M122").Value = -1035.5716
$ws.Range("N122").Value = -13918
$ws.Range("H132").Value = 2206.6538
$ws.Range("I132").Value = 2094.92
$ws.Range("K132").Value = 6284.76
$ws.Range("M132").Value = -3754.76
$ws.Range("H135").Value = 863.3333
$ws.Range("I135").Value = 863.3333
$ws.Range("K135").Value = 7769.9997
$ws.Range("M135").Value = -5234.9997
$ws.Range("H137").Value = 3968.3333
$ws.Range("I137").Value = 3949
$ws.Range("K137").Value = 11847
$ws.Range("M137").Value = -9297

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3542.8
$ws.Range("I2").Value = 3400.3333
$ws.Range("K2").Value = 3400.3333
$ws.Range("M2").Value = -3287.3333
$ws.Range("H4").Value = 449
$ws.Range("J4").Value = 210.66667
$ws.Range("L4").Value = 210.66667
$ws.Range("N4").Value = -442.66667
$ws.Range("H45").Value = 1798.4
$ws.Range("I45").Value = 996
$ws.Range("K45").Value = 996
$ws.Range("M45").Value = -619
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H110").Value = 1223
$ws.Range("I110").Value = 1223
$ws.Range("K110").Value = 1223
$ws.Range("M110").Value = 822
$ws.Range("H116").Value = 3542.8
$ws.Range("I116").Value = 3400.3333
$ws.Range("K116").Value = 3400.3333
$ws.Range("M116").Value = -1106.3333
$ws.Range("H132").Value = 1002.4375
$ws.Range("I132").Value = 1002.4375
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3007.3125
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -477.3125
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3542.8
$ws.Range("I3").Value = 3400.3333
$ws.Range("K3").Value = 3400.3333
$ws.Range("M3").Value = -3286.3333
$ws.Range("H54").Value = 28000
$ws.Range("I54").Value = 28000
$ws.Range("K54").Value = 28000
$ws.Range("M54").Value = -27516
$ws.Range("H134").Value = 3649.3157
$ws.Range("I134").Value = 3649.3157
$ws.Range("K134").Value = 10947.9471
$ws.Range("M134").Value = -8412.947100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1088.3636
$ws.Range("J16").Value = 595
$ws.Range("L16").Value = 595
$ws.Range("N16").Value = -1169
$ws.Range("H31").Value = 3209.875
$ws.Range("J31").Value = 4162.75
$ws.Range("L31").Value = 4162.75
$ws.Range("N31").Value = -4752.75
$ws.Range("H34").Value = 3209.875
$ws.Range("J34").Value = 4162.75
$ws.Range("L34").Value = 4162.75
$ws.Range("N34").Value = -4566.75
$ws.Range("H58").Value = 5069.7646
$ws.Range("I58").Value = 4042.6428
$ws.Range("K58").Value = 4042.6428
$ws.Range("M58").Value = -3839.6428
$ws.Range("H99").Value = 4139.25
$ws.Range("I99").Value = 4879
$ws.Range("J99").Value = 3695.4
$ws.Range("K99").Value = 4879
$ws.Range("L99").Value = 3695.4
$ws.Range("M99").Value = -3381
$ws.Range("N99").Value = -6691.4
$ws.Range("H107").Value = 769.5
$ws.Range("I107").Value = 508.27274
$ws.Range("J107").Value = 1180
$ws.Range("K107").Value = 508.27274
$ws.Range("L107").Value = 1180
$ws.Range("M107").Value = 1411.72726
$ws.Range("N107").Value = -5020
$ws.Range("H109").Value = 49995
$ws.Range("J109").Value = 49995
$ws.Range("L109").Value = 49995
$ws.Range("N109").Value = -52075
$ws.Range("H113").Value = 1088.3636
$ws.Range("J113").Value = 595
$ws.Range("L113").Value = 595
$ws.Range("N113").Value = -4935
$ws.Range("H126").Value = 4139.25
$ws.Range("I126").Value = 4879
$ws.Range("J126").Value = 3695.4
$ws.Range("K126").Value = 14637
$ws.Range("L126").Value = 11086.2
$ws.Range("M126").Value = -12167
$ws.Range("N126").Value = -16026.2
$ws.Range("H136").Value = 5069.7646
$ws.Range("I136").Value = 4042.6428
$ws.Range("K136").Value = 12127.9284
$ws.Range("M136").Value = -9577.928400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 1780.1666
$ws.Range("I46").Value = 69
$ws.Range("J46").Value = 2122.4
$ws.Range("K46").Value = 207
$ws.Range("L46").Value = 6367.200000000001
$ws.Range("M46").Value = -116
$ws.Range("N46").Value = -6549.200000000001
$ws.Range("H108").Value = 589.8333
$ws.Range("I108").Value = 589.8333
$ws.Range("K108").Value = 1769.4999
$ws.Range("M108").Value = 1110.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 70
$ws.Range("I2").Value = 87.22221999999999
$ws.Range("J2").Value = 18.333334
$ws.Range("K2").Value = 87.22221999999999
$ws.Range("L2").Value = 18.333334
$ws.Range("M2").Value = 25.77778000000001
$ws.Range("N2").Value = -244.333334
$ws.Range("H80").Value = 2996.5
$ws.Range("I80").Value = 2005
$ws.Range("K80").Value = 2005
$ws.Range("M80").Value = -1007
$ws.Range("H83").Value = 2996.5
$ws.Range("I83").Value = 2005
$ws.Range("K83").Value = 10025
$ws.Range("M83").Value = -5033
$ws.Range("H104").Value = 29890.334
$ws.Range("J104").Value = 29890.334
$ws.Range("L104").Value = 29890.334
$ws.Range("N104").Value = -36878.334
$ws.Range("H113").Value = 2336.6667
$ws.Range("I113").Value = 2336.6667
$ws.Range("K113").Value = 2336.6667
$ws.Range("M113").Value = -166.6667000000002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 68620.5
$ws.Range("I7").Value = 68620.5
$ws.Range("K7").Value = 68620.5
$ws.Range("M7").Value = -68508.5
$ws.Range("H22").Value = 700
$ws.Range("I22").Value = 700
$ws.Range("K22").Value = 700
$ws.Range("M22").Value = -405
$ws.Range("H27").Value = 700
$ws.Range("I27").Value = 700
$ws.Range("K27").Value = 700
$ws.Range("M27").Value = -593
$ws.Range("H46").Value = 4100
$ws.Range("I46").Value = 2000
$ws.Range("J46").Value = 4333.3335
$ws.Range("K46").Value = 2000
$ws.Range("L46").Value = 4333.3335
$ws.Range("M46").Value = -1812
$ws.Range("N46").Value = -4709.3335
$ws.Range("H55").Value = 245.83333
$ws.Range("I55").Value = 164.66667
$ws.Range("J55").Value = 327
$ws.Range("K55").Value = 164.66667
$ws.Range("L55").Value = 327
$ws.Range("M55").Value = 8.333329999999989
$ws.Range("N55").Value = -673
$ws.Range("H126").Value = 68620.5
$ws.Range("I126").Value = 68620.5
$ws.Range("K126").Value = 205861.5
$ws.Range("M126").Value = -203391.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 25000
$ws.Range("J104").Value = 25000
$ws.Range("L104").Value = 25000
$ws.Range("N104").Value = -31988

